# Generate Report for Handoff
#
# - Status moves from "In Translation" to "Ready for handoff" for both
#   locales (reflected on the Overview sheet as well as each locale sheet).
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#   timestamps are refreshed to the new handoff-generation time.
# - The Status / date columns are widened slightly to fit the new text.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---
$ovw.Range("E2").Value = "Ready for handoff"
$ovw.Range("F2").Value = "Ready for handoff"
$zh.Range("C2").Value  = "Ready for handoff"
$de.Range("C2").Value  = "Ready for handoff"

# --- Refreshed handoff timestamps ---
$ovw.Range("G2").Value = "2016-09-04 15:02:28"
$de.Range("H2").Value  = "2016-09-04 15:02:28"
$zh.Range("H2").Value  = "2016-09-04 15:02:24"

# --- Widen the Status/date columns to fit the new text ---
$newWidth = 16.333333333333332
$ovw.Columns("E:F").ColumnWidth = $newWidth
$zh.Columns("C:C").ColumnWidth  = $newWidth
$de.Columns("C:C").ColumnWidth  = $newWidth
